$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''57.049.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = '''2.320.85'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -1.96%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.44%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''532.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +2.01%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''132.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -2.62%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D8').Value = '''0.536'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.91%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''2.344.64'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -1.93%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = '''  -1.39%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = '''  -0.04%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''5.31'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -2.63%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''0.346'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +0.60%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''2.737.72'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -1.96%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''23.47'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -4.21%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''57.103.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -0.71%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '''  -2.28%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''2.364.47'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -0.61%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''337.98'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +2.00%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''10.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -1.96%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''6.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +2.47%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''4.17'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -1.90%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '''  +0.01%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''62.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +1.20%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').Value = '''InternetComputer(DFINITY)'
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = '''8.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  +0.40%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value = '''Kaspa'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = '''0.166'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  +0.24%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''0.996'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +0.13%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''1.35'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -0.75%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''173.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +2.98%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''1.73'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.93%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''0.0₃0726'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -3.01%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = '''  -3.32%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''18.55'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -0.41%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '''  -0.03%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''0.991'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -0.21%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = '''  -4.45%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '''  -0.71%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''4.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -1.47%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''39.25'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +1.38%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = '''  -2.98%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''149.22'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -1.23%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = '''RenderToken'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = '''5.41'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +0.64%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = '''PolygonEcosystemToken'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = '''https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = '''0.375'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -3.44%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = '''Filecoin'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = '''3.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -1.90%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = '''Bittensor'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = '''281.88'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -1.03%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''0.0931'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -1.24%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''0.0502'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -1.81%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''18.88'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +3.07%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '''  -1.38%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  -1.81%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = '''  +5.11%  '
$ws.Range('E51').Style = 'Normal'
